$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Subgroups")

$ws.Range("A8").Value  = "{'Hobby': np.int64(1), 'Student': np.int64(1)}"
$ws.Range("A9").Value  = "{'Gender': np.int64(1), 'Hobby': np.int64(1)}"
$ws.Range("A10").Value = "{'Hobby': np.int64(1), 'SexualOrientation': np.int64(1)}"
$ws.Range("A12").Value = "{'Gender': np.int64(1), 'Student': np.int64(1)}"
$ws.Range("A14").Value = "{'HDI': np.int64(1), 'Student': np.int64(1)}"
$ws.Range("A15").Value = "{'Gender': np.int64(1), 'SexualOrientation': np.int64(1)}"
$ws.Range("A16").Value = "{'Gender': np.int64(1), 'HDI': np.int64(1)}"
$ws.Range("A17").Value = "{'HDI': np.int64(1), 'SexualOrientation': np.int64(1)}"
$ws.Range("A18").Value = "{'Gender': np.int64(1), 'Hobby': np.int64(1), 'SexualOrientation': np.int64(1)}"
$ws.Range("A19").Value = "{'Gender': np.int64(1), 'Student': np.int64(1), 'SexualOrientation': np.int64(1)}"
$ws.Range("A20").Value = "{'Gender': np.int64(1), 'HDI': np.int64(1), 'Student': np.int64(1)}"
$ws.Range("A21").Value = "{'HDI': np.int64(1), 'Student': np.int64(1), 'SexualOrientation': np.int64(1)}"
$ws.Range("A22").Value = "{'Gender': np.int64(1), 'HDI': np.int64(1), 'SexualOrientation': np.int64(1)}"
